$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.683.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.626.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.625.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("E14").Value = "  +5.78%  "
$ws.Range("E15").Value = "  +5.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.644.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.629.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "384.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.69%  "
$ws.Range("E20").Value = "  +6.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.763.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("E30").Value = "  +7.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "530.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("E32").Value = "  +3.66%  "
$ws.Range("E33").Value = "  +7.31%  "
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +6.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.110"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.62%  "
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("E49").Value = "  +5.57%  "
$ws.Range("E50").Value = "  +8.22%  "
$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").Value = "  +4.50%  "
